# Applies the "01st-24 report" update to the Route Cost RSO workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Route")

# Update the report date in L3 (merged L3:M3) from text "30/9/2024" to an
# actual date value (1/10/2024), keeping the existing date number format.
$ws.Range("L3").Value = Get-Date -Year 2024 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Update the route-cost amounts for RSO 02, RSO 04 and RSO 03 rows.
# Dependent formulas (I, J, K, L columns and the totals row) recalculate
# automatically.
$ws.Range("D7").Value = 150
$ws.Range("D9").Value = 150
$ws.Range("D10").Value = 200

$excel.CalculateFullRebuild()
